$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "~TFM_UPD" table block (rows 13-15), mirroring the existing
# "~TFM_INS" block (rows 9-11), to make combined electricity consumption of
# DCs on ANNUAL level equal to the end-use demand of DCs ---

# Row 13: table-name marker cell, same formatting as B9 ("~TFM_INS")
$ws.Range("B9").Copy($ws.Range("B13")) | Out-Null
$ws.Range("B13").Value = "~TFM_UPD"

# Row 14: header row, identical to the existing header row 10
$ws.Range("B10:I10").Copy($ws.Range("B14:I14")) | Out-Null

# Row 15: data row, styled like row 11
$ws.Range("D11").Copy($ws.Range("D15")) | Out-Null
$ws.Range("D15").Value = "EFF"

$ws.Range("F11").Copy($ws.Range("F15")) | Out-Null
$ws.Range("F15").Value = 1.073

$ws.Range("G11").Copy($ws.Range("G15")) | Out-Null
$ws.Range("G15").Formula = "=F15"

$ws.Range("H11").Copy($ws.Range("H15")) | Out-Null
$ws.Range("H15").Value = "SRVDCE-CS"

# Threaded comment on F15 explaining the rationale
$excel.UserName = "Balyk, Olexandr"
$comment = $ws.Range("F15").AddCommentThreaded("lower consumption electricity that is not used for cooling. I.e. el for cooling + el for IT ca. = dc end-use demand")

# Update the active selection to reflect where editing ended
$ws.Range("I15").Select() | Out-Null
